$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 9253
$ws.Range("B2").Value = "Cauã Farias"
$ws.Range("C2").Value = "Marketing"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45094
$ws.Range("G2").Value = 11883.42

# Row 3
$ws.Range("A3").Value = 50180
$ws.Range("B3").Value = "Sr. João Vitor Castro"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45105
$ws.Range("G3").Value = 8534.2

# Row 4
$ws.Range("A4").Value = 83379
$ws.Range("B4").Value = "Rafaela Sales"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45092
$ws.Range("G4").Value = 10488.67

# Row 5
$ws.Range("A5").Value = 53694
$ws.Range("B5").Value = "Paulo Lopes"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("F5").Value = 45102
$ws.Range("G5").Value = 11993.64

# Row 6
$ws.Range("A6").Value = 46615
$ws.Range("B6").Value = "Ana Sophia Azevedo"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45105
$ws.Range("G6").Value = 7932.62

# Row 7
$ws.Range("A7").Value = 93351
$ws.Range("B7").Value = "Pedro Lucas Melo"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45098
$ws.Range("G7").Value = 6814.74

# Row 8
$ws.Range("A8").Value = 82498
$ws.Range("B8").Value = "Dr. Fernando Vieira"
$ws.Range("C8").Value = "Financeiro"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45083
$ws.Range("G8").Value = 4081.33

# Row 9
$ws.Range("A9").Value = 32687
$ws.Range("B9").Value = "Otávio Teixeira"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45092
$ws.Range("G9").Value = 9583.42

# Row 10
$ws.Range("A10").Value = 55810
$ws.Range("B10").Value = "Sra. Agatha Moreira"
$ws.Range("C10").Value = "Marketing"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("F10").Value = 45082
$ws.Range("G10").Value = 9602.25

# Row 11
$ws.Range("A11").Value = 44745
$ws.Range("B11").Value = "Luiz Gustavo Nascimento"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45096
$ws.Range("G11").Value = 5669.3
